$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed/Modified) date column for rows 2-5
# from serial date 45212 (2023-10-13) to 45221 (2023-10-22),
# preserving existing number formatting/style on the cells.
$ws.Range("C2:C5").Value = 45221
